$d = $word.ActiveDocument

function Insert-XmlAtRange($rng, $bodyXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# --- 1. "Note:" paragraph: tweak last run and append a new run about keystrokes ---
$p = $d.Paragraphs(7)
$rng = $p.Range
$xml = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Note</w:t></w:r>' +
       '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">: The basic layout of the document was decided right away and wasn&apos;t </w:t><w:tab/>' +
       '<w:t xml:space="preserve">changed. This included a menu at the top with your standard File and Edit </w:t><w:tab/>' +
       '<w:t xml:space="preserve">items, a tool panel right underneath, a text area for the document itself, and </w:t><w:tab/>' +
       '<w:t xml:space="preserve">a panel at the bottom (which in theory would include information about the </w:t><w:tab/>' +
       '<w:t>state of the document</w:t></w:r>' +
       '<w:r><w:rPr></w:rPr><w:t>)</w:t></w:r>' +
       '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">. All of these are the basic elements of a common text </w:t><w:tab/>' +
       '<w:t xml:space="preserve">editor, which we stuck with for the sake of learnability. </w:t></w:r>' +
       '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">Also, the standard </w:t><w:tab/>' +
       '<w:t xml:space="preserve">keystrokes for undo/redo as well as other functionalities have been </w:t><w:tab/>' +
       '<w:t>implemented.</w:t></w:r>' +
       '</w:p>'
Insert-XmlAtRange $rng $xml

# --- 2. Typo fix: "somtimes" -> "sometimes" (rebuild paragraph to keep <w:tab/> intact) ---
$p19 = $d.Paragraphs(19)
$rng19 = $p19.Range
$xml19 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
         '<w:t xml:space="preserve">-Moving cursor: Also technically implemented, sometimes moving the cursor will </w:t><w:tab/>' +
         '<w:t>cause the history to go out of sync, giving undesirable behavior</w:t></w:r></w:p>'
Insert-XmlAtRange $rng19 $xml19

# --- 3. Last bugs paragraph: extend sentence, then append Scenarios section ---
# This is the last paragraph in the document, so its Range.End sits at the very
# end of the body; excluding that final paragraph-mark position keeps InsertXML
# from leaving a stray empty paragraph behind.
$p2 = $d.Paragraphs(20)
$rng2raw = $p2.Range
$rng2 = $d.Range($rng2raw.Start, $rng2raw.End - 1)
$xml2 =
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
  '<w:t>These cause of these bugs are not well known, making it difficult to fix. None of these bugs crash the program however, so no data is lost.</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Scenarios</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Rob, an inexperienced user</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
  '<w:t>Sally head from his friends Sally and John that there is this new amazing text editor called Extreme Editor that can be used. Rob, who has had some experience with other text editors, gives it a try. He decides to write a poem in it. He is familiar with the basic undo/redo buttons in other programs, and used it a lot while making the poem. Rob, though not as technical as his friends Sally and John, is very good about saving his work. When asked if he wanted to save his undo/redo history, he selects yes. By some coincidence, right after this, his computer crashed. Rob was pleased to find that when reopening his file in Extreme Editor, his history was intact and he could undo/redo his past work.</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:i/><w:i/><w:iCs/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Sally, an occasional user</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr></w:rPr></w:pPr>' +
  '<w:r><w:rPr></w:rPr><w:tab/><w:t xml:space="preserve">Sally </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">uses </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">Extreme Editor </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>occasionally</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">, and is aware of the shortcuts and features a new user may not </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>know about</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">. She is going to write a </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">Java </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">program to </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>manipulate large amounts of data</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve">. When she wants to change a section, she opens the undo history window and selects the individual items she wants changed. This allows her to have great control over the modifications, while saving time compared to manually selecting the text and erasing/re-typing it. </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>She can easily r</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>emov</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>e</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve"> entire </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>p</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>aragraphs, and replac</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>e</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t xml:space="preserve"> them when needed, at a few mouse presses. Sally has written code in her main </w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>method</w:t></w:r>' +
  '<w:r><w:rPr></w:rPr><w:t>, but wants to move it to a separate method while also removing certain lines. She selects the relevant section, but instead of copying or moving, she accidently deletes it. She opens the history menu and undo&#8217;s the lines she wants to keep, but leaving the lines to be erased, erased.</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>John, an experienced user</w:t></w:r></w:p>' +

  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="480"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/>' +
  '<w:t>John has been tasked by his employer to store important employee information in a text document. John, an everyday user of Extreme Editor, begins his task. Halfway through, he recieves a call from his boss asking him to exclude a addresses from the document. John was quick to open the undo/redo history, and promptly selects all address insertions. Wanting to make sure the result was as he wanted, he decides to preview his changes. Satisfied with the results, John applies the changes, and reports to his boss that the change has been made. His boss, impressed at his efficiency, gives him a promotion.</w:t></w:r></w:p>'

Insert-XmlAtRange $rng2 $xml2

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
